$wb = $excel.ActiveWorkbook

# --- Fix "m2pc BAU-S - starting point": rows 4:9 (years 2025-2050) were
# cascading the 2020 value down via shared formulas (=C3, =D3, ...).
# Point them instead at the matching row of "m2pc S - starting point".
$wsBauSStart = $wb.Worksheets.Item("m2pc BAU-S - starting point")
$wsBauSStart.Range("C4:G9").Formula = "='m2pc S - starting point'!C4"

# --- Update view state (zoom / selection) on the sheets that were touched,
# finishing on "m2pc BAU-S - starting point" so it ends up the active tab.

$wsSStart = $wb.Worksheets.Item("m2pc S - starting point")
$wsSStart.Activate()
$excel.ActiveWindow.Zoom = 280
$wsSStart.Range("C4:G9").Select()

$ws0 = $wb.Worksheets.Item("m2pc 0")
$ws0.Activate()
$ws0.Range("J1:J1048576").Select()

$wsBau = $wb.Worksheets.Item("m2pc BAU")
$wsBau.Activate()
$wsBau.Range("C10").Select()

$wsS = $wb.Worksheets.Item("m2pc S")
$wsS.Activate()
$wsS.Range("J1:J1048576").Select()

$wsBauS = $wb.Worksheets.Item("m2pc BAU-S")
$wsBauS.Activate()
$wsBauS.Range("D5").Select()

$wsBauSStart.Activate()
$excel.ActiveWindow.Zoom = 280
$wsBauSStart.Range("B13").Select()
